$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B6 value from "Baby Boomers" to "Boomers"
$ws.Range("B6").Value = "Boomers"

# Update the selection shown in the sheet view to B6
$ws.Activate()
$ws.Range("B6").Select()
